# "Add files via upload" — a re-upload of the daily error-count tracker
# with one more day of data (10/20) appended to the table, and the
# active-cell selection left wherever the author had last clicked.
#
# (The surrounding cosmetic diffs in the canonical XML — workbookView
# window geometry, x15ac:absPath, xr:revisionPtr GUIDs, sheetFormatPr
# defaultRowHeight/x14ac:dyDescent, and the sub-pixel bestFit column
# width — are artifacts Excel stamps from the local machine/session that
# saved the file, not values exposed on the Excel object model for a
# script to set; they are left untouched here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data for 2025-10-20, appended right after the existing
# 2025-10-17 row (row 9) and before the trailing blank template rows.
$ws.Range("A10").Value = 45950
$ws.Range("B10").Value = 541
$ws.Range("C10").Value = 11
$ws.Range("D10").Value = 0.020183486238532111
$ws.Range("E10").Value = 11
$ws.Range("F10").Value = 0.97981651376146783

# Leave the selection where the author's cursor ended up.
$ws.Range("K9").Select()
